# Logboek Eilish Van Der Snickt - "Instellingenscherm en functionaliteiten toevoegen"
# Adds the Week-5 total to row 15 (B15), inserts a new "Week 5" header row (17)
# and a new logboek entry row (18) with date, duration, description and a
# youtube link, mirroring the existing rows' layout/styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 15: add the week-5 total-time cell (B15), matching the style used
#     for the other "Week N" total cells (e.g. B10).
$ws.Range("B15").Value2 = "2 uur 15 minuten"
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B15").PasteSpecial($xlPasteFormats) | Out-Null

# --- New row 18: the logboek entry for the new week.
#     Fill Q18 first, then A17, then C18 so new shared strings are created
#     in the same order as the source edit.
$ws.Range("Q18").Value2 = "https://www.youtube.com/watch?v=eX-TdY6bLdg"

# --- Row 17: turn the (until now) mostly-empty row into the "Week 5" header row.
$ws.Range("A17").Value2 = "Week 5"
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A17").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B10").Copy() | Out-Null
$ws.Range("B17").PasteSpecial($xlPasteFormats) | Out-Null

# --- Row 18 continued: date, duration and description.
$ws.Range("C18").Value2 = "Instellingen scherm verder uitwerken, navigatie aanpassen, popupscherm maken, gebruikers verwijderen en beginnende code voor wachtwoord wijziging"

$ws.Range("A18").Value2 = 43540
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A18").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B18").Value2 = "4 uur"

$excel.CutCopyMode = 0

# --- Update the remembered selection to match the saved workbook.
$ws.Range("L23").Select() | Out-Null
